# Fill in the "Day 17" colony counts (column W) for data rows 3-34.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dayCounts = [ordered]@{
    3  = 69
    4  = 103
    5  = 95
    6  = 101
    7  = 147
    8  = 55
    9  = 101
    10 = 81
    11 = 65
    12 = 32
    13 = 59
    14 = 66
    15 = 18
    16 = 62
    17 = 112
    18 = 30
    19 = 84
    20 = 55
    21 = 71
    22 = 56
    23 = 165
    24 = 183
    25 = 93
    26 = 95
    27 = 54
    28 = 42
    29 = 59
    30 = 48
    31 = 184
    32 = 138
    33 = 138
    34 = 60
}

foreach ($row in $dayCounts.Keys) {
    $ws.Range("W$row").Value = $dayCounts[$row]
}
